$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of resale-number data for 2024-01-09 14:40:26 (row 38).
$rowIndex = 38

# Text-like columns (A-D) must stay literal text, not get auto-converted
# by Excel's smart typing (dates/leading-zero numbers). Force the cell to
# Text format before assignment, then clear the explicit formatting back
# off so no stray style id is left behind on the cell.
$textValues = @{
    1 = "2024-01-09"   # A - Date
    2 = "14:40:26"     # B - Time
    3 = "Tuesday"      # C - Weekday
    4 = "01"           # D - Week
}

foreach ($col in $textValues.Keys) {
    $cell = $ws.Cells.Item($rowIndex, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $textValues[$col]
    $cell.ClearFormats()
}

# Numeric columns (E-T).
$numericValues = @{
    5  = 139476   # E - Beijing
    6  = 142694   # F - Guangzhou
    7  = 171996   # G - Suzhou
    8  = 147537   # H - Hangzhou
    9  = -1       # I - Nanjing
    10 = 118186   # J - Xi_an
    11 = 224707   # K - Chengdu
    12 = 250270   # L - Chongqing
    13 = 185124   # M - Tianjin
    14 = 110384   # N - Hefei
    15 = 40657    # O - Fuzhou
    16 = 30845    # P - Xiamen
    17 = 72562    # Q - Changsha
    18 = -1       # R - Shanghai
    19 = 41849    # S - Shenzhen
    20 = -1       # T - Wuhan
}

foreach ($col in $numericValues.Keys) {
    $ws.Cells.Item($rowIndex, $col).Value = $numericValues[$col]
}
